$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = '#caezar'
$ws.Range("C2").Value = 'Caezar'
$ws.Range("D2").ClearContents()

$ws.Range("B3").Value = '#philippus'
$ws.Range("C3").Value = 'Philippus'
$ws.Range("D3").ClearContents()

$ws.Range("B4").Value = '#prolomeus'
$ws.Range("C4").Value = 'Prolomeus'
$ws.Range("D4").ClearContents()

$ws.Range("B5").Value = '#ptolomeus'
$ws.Range("C5").Value = 'Ptolomeus'
$ws.Range("D5").ClearContents()

$ws.Range("B6").Value = '#charmione'
$ws.Range("C6").Value = 'Charmione'
$ws.Range("D6").ClearContents()

$ws.Range("B7").Value = '#photinus'
$ws.Range("C7").Value = 'Photinus'
$ws.Range("D7").ClearContents()

$ws.Range("B8").Value = '#rome'
$ws.Range("C8").Value = 'Rome'
$ws.Range("D8").ClearContents()

$ws.Range("B9").Value = '#septimius'
$ws.Range("C9").Value = 'Septimius'
$ws.Range("D9").ClearContents()

$ws.Range("B10").Value = '#ornelia'
$ws.Range("C10").Value = 'Ornelia'
$ws.Range("D10").ClearContents()

$ws.Range("B11").Value = '#cornelia'
$ws.Range("C11").Value = 'Cornelia'
$ws.Range("D11").ClearContents()

$ws.Range("B12").Value = '#kornelia'
$ws.Range("C12").Value = 'Kornelia'
$ws.Range("D12").ClearContents()

$ws.Range("B13").Value = '#anthonius'
$ws.Range("C13").Value = 'Anthonius'
$ws.Range("D13").ClearContents()

$ws.Range("B14").Value = '#achoreus'
$ws.Range("C14").Value = 'Achoreus'
$ws.Range("D14").ClearContents()

$ws.Range("B15").Value = '#cleopatra,'
$ws.Range("C15").Value = 'Cleopatra,'
$ws.Range("D15").ClearContents()

$ws.Range("B16").Value = '#aegypte'
$ws.Range("C16").Value = 'AEgypte'
$ws.Range("D16").ClearContents()

$ws.Range("B17").Value = '#achillas'
$ws.Range("C17").Value = 'Achillas'
$ws.Range("D17").ClearContents()

$ws.Range("B18").Value = '#kornelia:'
$ws.Range("C18").Value = 'Kornelia:'

$ws.Range("B19").Value = '#wysheid'
$ws.Range("C19").Value = 'Wysheid'

$ws.Range("B20").Value = '#lepidus'
$ws.Range("C20").Value = 'Lepidus'

$ws.Range("B21").Value = '#petolomeus'
$ws.Range("C21").Value = 'Petolomeus'

$ws.Range("B22").Value = '#cleopatra'
$ws.Range("C22").Value = 'Cleopatra'

$ws.Range("B23").Value = '#staatkunde'
$ws.Range("C23").Value = 'Staatkunde'
